# Rename the sheet from "18_0s" to "population"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "population"

# Update the saved view state: clear the scrolled-to-row-71 view and
# move the selection/active cell to C6 (instead of the full-column
# A:A selection), matching the author's latest interactive state.
$ws.Activate()
$ws.Range("C6").Select()
